$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# and they would lose their original text representation (trailing zeros, etc).
$ws.Range("D2").Value = "59.293.01"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "2.509.20"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.05"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.69"
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").Value = "2.511.00"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.39"
$ws.Range("E12").Value = "  -3.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("E13").Value = "  -2.74%  "
$ws.Range("D14").Value = "2.959.46"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.23"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "59.175.54"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "2.508.13"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.75"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.46"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("E25").Value = "  -4.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  -3.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.76"
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("D30").Value = "0.0₃0775"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.64"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("E35").Value = "  -9.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.51"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.19"
$ws.Range("E37").Value = "  -5.90%  "
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.67"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("E42").Value = "  -8.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "279.94"
$ws.Range("E43").Value = "  -7.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0936"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.63"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0223"
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.73"
$ws.Range("E51").Value = "  -3.28%  "
